$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102: Epidemiologic term (Data Content Type)
# Set cell values in the same order the original commit's shared-string
# table grew (definition, term, ... ) so sharedStrings indices line up.
$ws.Range("C102").Value = "Relating to the study of the distribution and determinants of health-related states or events (including disease) in populations, and the application of this study to the control of diseases and other health problems."
$ws.Range("B102").Value = "Epidemiologic"
$ws.Range("A102").Value = "Data Content Type"
$ws.Rows.Item(102).RowHeight = 29

# Row 103: Proband term (Data Element)
$ws.Range("B103").Value = "Proband"
$ws.Range("A103").Value = "Data Element"
$ws.Range("C103").Value = "A proband is a person in a family to receive genetic counseling and/or testing for a suspected hereditary risk or diagnosed disease. A proband may or may not be affected with the disease in question. If the value is true, then the case subject may have been diagnosed with the disease under studied. If the value is false, then the case subject is a member of the family of a proband study participant. The proband indicator for the case carries over to a sample taken from a case subject."
$ws.Range("D103").Value = "NCI Dictionary"
$ws.Rows.Item(103).RowHeight = 72.5

# Reference URL for the Proband definition, as a real hyperlink.
$ws.Hyperlinks.Add($ws.Range("E103"), "https://www.cancer.gov/publications/dictionaries/genetics-dictionary/def/proband")

# Adjust column widths to match the refreshed glossary layout.
$ws.Columns.Item(1).ColumnWidth = 36.36328125
$ws.Columns.Item(2).ColumnWidth = 19.6328125
$ws.Columns.Item(4).ColumnWidth = 12.90625

# Restore the view to show the newly-added rows, as the author left it.
$ws.Application.ActiveWindow.ScrollRow = 94
$ws.Range("B106").Select()
